$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", from, locale, timezone)+" - "+dateTool.format("YYYY-MM-dd HH:mm:ss", to, locale, timezone)}'
$ws.Range("B9").Value = '${dateTool.format("YYYY-MM-dd HH:mm:ss", position.fixTime, locale, timezone)}'

$ws.Range("B2").Select()
